# Response-to-reviews.docx edits
# - strip leading "XXX " from several paragraphs now that they have been handled
# - reflow "In general, we have ..." -> "We have in general ..." (with odd run splits
#   matching how the author's word processor recorded the retyped text)
# - merge two runs that made up "... is "Never use the passive ..." ... phrasing."
#   into a single run
# - append a new "XXX or have we?" aside to the "We have retained ..." paragraph
# - drop the leading "XXX " from the Tschopp Method 4 paragraph
# - flip Normal style's overflow punctuation flag off
# - mint a new (empty) ListLabel9 character style

$d = $word.ActiveDocument

# 1) Drop "XXX " prefixes on four paragraphs.
$d.Content.Find.Execute(
    "XXX Marek" + [char]8217 + "s multiple comments",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Marek" + [char]8217 + "s multiple comments", 2) | Out-Null

$d.Content.Find.Execute(
    "XXX We also accept his point",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We also accept his point", 2) | Out-Null

$d.Content.Find.Execute(
    "XXX As recommended, we have greatly reduced",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "As recommended, we have greatly reduced", 2) | Out-Null

$d.Content.Find.Execute(
    "XXX Tschopp asks why",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Tschopp asks why", 2) | Out-Null

# 2) Re-work the "In general, we have followed ..." sentence into
#    "We have in general followed ...", splitting it into the same
#    four runs the author ended up with, and merge the trailing
#    "and note that the fourth of " run into the last of those four.
$rng = $d.Content
$rng.Find.Execute("In general, we have followed the specific comments attached to Marek" + [char]8217 + "s review, but with some exceptions. In particular, we do not agree with the suggestion that expunging the pronoun " + [char]8220 + "we" + [char]8221 + " throughout and substituting passive voice would improve the manuscript, ") | Out-Null
$target = $rng.Duplicate
$target.MoveEnd(1, ("and note that the fourth of ").Length) | Out-Null
$target.Text = "We" + "e have " + "i" + "n general followed the specific comments attached to Marek" + [char]8217 + "s review, but with some exceptions. In particular, we do not agree with the suggestion that expunging the pronoun " + [char]8220 + "we" + [char]8221 + " throughout and substituting passive voice would improve the manuscript, and note that the fourth of "

# 3) Merge the "is ... active"" / ". We have changed ..." runs (simple concatenation,
#    no wording change) by replacing the straight-quote-delimited phrase through to
#    "active phrasing." in one go, preserving the existing xml:space="preserve" run.
$d.Content.Find.Execute(
    [char]8220 + "Never use the passive where you can use the active" + [char]8221 + [char]10 + ". We have changed",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# 4) Append "XXX or have we?" after the open peer-review paragraph, putting a
#    trailing space on the existing sentence and a new run for the aside.
$rng2 = $d.Content
$rng2.Find.Execute("We have retained the section on open peer-review, contrary to Marek" + [char]8217 + "s recommendation, as it is important to us and relevant to the origin of the present paper" + [char]8217 + "s core question.") | Out-Null
$rng2.InsertAfter(" XXX or have we?")

# 5) Normal style: stop allowing punctuation to overflow the margin.
$d.Styles("Normal").ParagraphFormat.HangingPunctuation = $false

# 6) Mint a new, empty "ListLabel 9" character style (id ListLabel9).
$newStyle = $d.Styles.Add("ListLabel9", 2)
$newStyle.NameLocal = "ListLabel 9"
$newStyle.QuickStyle = $true

Write-Output "done"
